# Automating the Input files for Transavia
# Updates computed GHI/DNI/DHI weather values on the "Daily" and "Hourly" sheets.

$wb = $excel.ActiveWorkbook

# --- Daily sheet (row 2) ---
$daily = $wb.Worksheets.Item("Daily")
$daily.Range("G2").Value = 4944.05
$daily.Range("H2").Value = 8205.9
$daily.Range("I2").Value = 994.36
$daily.Range("J2").Value = 3661.06
$daily.Range("K2").Value = 3765.83
$daily.Range("L2").Value = 1651.25

# --- Hourly sheet (rows 9-20) ---
$hourly = $wb.Worksheets.Item("Hourly")

# Row 9
$hourly.Range("K9").Value = 99.04000000000001
$hourly.Range("L9").Value = 172.18
$hourly.Range("M9").Value = 61.71

# Row 10
$hourly.Range("I10").Value = 670.8
$hourly.Range("J10").Value = 77.31
$hourly.Range("K10").Value = 292.39
$hourly.Range("L10").Value = 620.62
$hourly.Range("M10").Value = 80.29000000000001

# Row 11
$hourly.Range("H11").Value = 455.91
$hourly.Range("I11").Value = 773.17
$hourly.Range("K11").Value = 453.04
$hourly.Range("L11").Value = 739.8200000000001
$hourly.Range("M11").Value = 98.97

# Row 12
$hourly.Range("I12").Value = 830.72
$hourly.Range("J12").Value = 102.46
$hourly.Range("K12").Value = 569.3099999999999
$hourly.Range("L12").Value = 766.55
$hourly.Range("M12").Value = 121.32

# Row 13
$hourly.Range("H13").Value = 661.92
$hourly.Range("I13").Value = 860.77
$hourly.Range("J13").Value = 107.96
$hourly.Range("K13").Value = 603.1
$hourly.Range("L13").Value = 621.54
$hourly.Range("M13").Value = 198.81

# Row 14
$hourly.Range("H14").Value = 688.9299999999999
$hourly.Range("I14").Value = 870
$hourly.Range("J14").Value = 109.7
$hourly.Range("K14").Value = 575
$hourly.Range("L14").Value = 422.15
$hourly.Range("M14").Value = 290.45

# Row 15
$hourly.Range("H15").Value = 660.12
$hourly.Range("I15").Value = 860.21
$hourly.Range("J15").Value = 107.82
$hourly.Range("K15").Value = 506.95
$hourly.Range("L15").Value = 304.91
$hourly.Range("M15").Value = 308.4

# Row 16
$hourly.Range("H16").Value = 578.22
$hourly.Range("I16").Value = 829.47
$hourly.Range("J16").Value = 102.18
$hourly.Range("K16").Value = 336.83
$hourly.Range("L16").Value = 118.04
$hourly.Range("M16").Value = 265.89

# Row 17
$hourly.Range("H17").Value = 451.29
$hourly.Range("I17").Value = 770.9
$hourly.Range("J17").Value = 92.23999999999999
$hourly.Range("K17").Value = 112.83
$hourly.Range("M17").Value = 112.83

# Row 18
$hourly.Range("I18").Value = 666.64
$hourly.Range("K18").Value = 73.39
$hourly.Range("M18").Value = 73.39

# Row 19
$hourly.Range("I19").Value = 460.64
$hourly.Range("J19").Value = 51.83
$hourly.Range("K19").Value = 32.21
$hourly.Range("M19").Value = 32.21

# Row 20
$hourly.Range("I20").Value = 67.06999999999999
$hourly.Range("K20").Value = 2.59
$hourly.Range("M20").Value = 2.59
